$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rows 76 and 77 had their match data (columns F:V) swapped ---
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row76 = @{}
$row77 = @{}
foreach ($col in $cols) {
    $row76[$col] = $ws.Range($col + "76").Value2
    $row77[$col] = $ws.Range($col + "77").Value2
}
foreach ($col in $cols) {
    $ws.Range($col + "76").Value = $row77[$col]
    $ws.Range($col + "77").Value = $row76[$col]
}

# --- Step 2: append new row 91 (new match added to the bottom of the sheet) ---
$ws.Range("A90:V90").Copy()
$ws.Range("A91:V91").PasteSpecial(-4122)

$ws.Range("A91").Value = 90
$ws.Range("B91").Value = "portugal"
$ws.Range("C91").Value = "liga-portugal"
$ws.Range("D91").Value = "2023-2024"
$ws.Range("E91").Value = 45236.88541666666
$ws.Range("F91").Value = "SC Farense"
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = "Arouca"
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2.17
$ws.Range("K91").Value = "30/10/2023 14:42"
$ws.Range("L91").Value = 2.41
$ws.Range("M91").Value = "06/11/2023 21:07"
$ws.Range("N91").Value = 3.39
$ws.Range("O91").Value = "30/10/2023 14:42"
$ws.Range("P91").Value = 3.36
$ws.Range("Q91").Value = "06/11/2023 20:48"
$ws.Range("R91").Value = 3.61
$ws.Range("S91").Value = "30/10/2023 14:42"
$ws.Range("T91").Value = 3.18
$ws.Range("U91").Value = "06/11/2023 21:07"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/sc-farense-arouca/S81wX3rO/"
